$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.110.46'
$ws.Range('E2').Value = '  -1.81%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.514.50'
$ws.Range('E3').Value = '  -2.80%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.10'
$ws.Range('E5').Value = '  -2.89%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.17'
$ws.Range('E6').Value = '  -5.21%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.512.63'
$ws.Range('E7').Value = '  -2.87%  '

$ws.Range('E8').Value = '  -0.14%  '

$ws.Range('E9').Value = '  -0.98%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('E10').Value = '  -0.35%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.55'
$ws.Range('E11').Value = '  +4.71%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.431'
$ws.Range('E12').Value = '  -1.94%  '

$ws.Range('E13').Value = '  -3.37%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.12'
$ws.Range('E14').Value = '  -3.25%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.103.25'
$ws.Range('E15').Value = '  -3.00%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.517.61'
$ws.Range('E16').Value = '  -3.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.114.34'
$ws.Range('E17').Value = '  -1.94%  '

$ws.Range('E18').Value = '  -0.19%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.53'
$ws.Range('E19').Value = '  -1.27%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.46'
$ws.Range('E20').Value = '  -2.55%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.99'
$ws.Range('E21').Value = '  -1.22%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '451.53'
$ws.Range('E22').Value = '  -1.87%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.629'
$ws.Range('E23').Value = '  -1.67%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.10'
$ws.Range('E24').Value = '  +0.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.646.88'
$ws.Range('E25').Value = '  -3.12%  '

$ws.Range('E26').Value = '  +0.01%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000125'
$ws.Range('E27').Value = '  -8.52%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.73'
$ws.Range('E28').Value = '  -4.95%  '

$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.01'
$ws.Range('E29').Value = '  -6.19%  '

$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.67'
$ws.Range('E30').Value = '  -2.51%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.52'
$ws.Range('E31').Value = '  -3.61%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.170'
$ws.Range('E32').Value = '  -3.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.07%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.71'
$ws.Range('E34').Value = '  -2.60%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.21'
$ws.Range('E35').Value = '  -6.13%  '

$ws.Range('E36').Value = '  -5.97%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.502.67'
$ws.Range('E37').Value = '  -2.79%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.04'
$ws.Range('E38').Value = '  -3.93%  '

$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.29'
$ws.Range('E40').Value = '  -3.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.13%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '177.13'
$ws.Range('E42').Value = '  +0.59%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0907'
$ws.Range('E43').Value = '  -2.16%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.46'
$ws.Range('E44').Value = '  -2.86%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '30.83'
$ws.Range('E45').Value = '  -2.97%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.900'
$ws.Range('E46').Value = '  -1.14%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.99'
$ws.Range('E47').Value = '  +1.55%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.31'
$ws.Range('E48').Value = '  -4.67%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.65'
$ws.Range('E49').Value = '  -1.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.52'
$ws.Range('E50').Value = '  -10.50%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  -1.60%  '
